# Applies the "Added A02c for Project Ranking and Selection" edit:
#  1. Rename worksheet "ProjExploration" -> "ProjSelection"
#  2. Add a new rubric block (rows 20-25) to that sheet for
#     "Project Ranking and Selection"
#  3. Update sheet view (selection / top-left cell) to point at the new block

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjExploration")

# --- 1. Rename the sheet -----------------------------------------------
$ws.Name = "ProjSelection"

# --- 2. Add the new rubric block ----------------------------------------
# Row 21 is the header row (Score | Project Ranking and Selection), styled
# like the existing header rows (e.g. row 3 / row 10).
$ws.Range("C3:D3").Copy($ws.Range("C21:D21"))
$ws.Range("D21").Value = "Project Ranking and Selection"

# Row 22 -> score 3 (styled like row 4 / row 11)
$ws.Range("C4:D4").Copy($ws.Range("C22:D22"))
$ws.Range("D22").Value = "On time; Wiki contains all required information;  Information presented gives a clear and highly detailed picture of why the team selected the project that they did; Rankings and selection are clearly and well supported by thoughtful, reflective and insightful rationales that draw directly on the gathered information; Writing is clear, concise, well organized, uses complete sentences and proper grammar; Wiki is correctly linked, neatly formatted and easy to read."

# Row 23 -> score 2 (styled like row 5 / row 12)
$ws.Range("C5:D5").Copy($ws.Range("C23:D23"))
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = "On time; Wiki contains all required information; Information presented indicates why the team selected the project that they did; Rankings and selection are supported by rationales that mention gathered information; Writing may have minor issues but does not distract from meaning or understanding; Reviews are correctly linked but may be poorly formatted or difficult to read."

# Row 24 -> score 1 (styled like row 6 / row 13)
$ws.Range("C6:D6").Copy($ws.Range("C24:D24"))
$ws.Range("D24").Value = "On time; Wiki may be missing required information; Information in the Wiki does not give an adequate picture of why the team selected the project they did; Rankings may be insufficiently connected to the rationale or rationale may not use or may not align with provided information; Writing, formatting or grammar may interfere with understanding; Reviews may be incorrectly linked, poorly formatted or difficult to read."

# Row 25 -> score 0 (styled like row 7 / row 14, re-uses the same "Late,
# missing..." text already used by those rows)
$ws.Range("C7:D7").Copy($ws.Range("C25:D25"))

# --- Row heights -----------------------------------------------------
$ws.Rows.Item(20).RowHeight = 17
$ws.Rows.Item(21).RowHeight = 17
$ws.Rows.Item(22).RowHeight = 153
$ws.Rows.Item(23).RowHeight = 116
$ws.Rows.Item(24).RowHeight = 132
$ws.Rows.Item(25).RowHeight = 35

# --- 3. Update the sheet view -------------------------------------------
$ws.Application.Goto($ws.Range("D25"))
$ws.Range("D25").Select()

Write-Host "Applied Project Ranking and Selection rubric block."
